# Simulated Wild Card round and logged it
# Update the "H" (Home) row totals on both the OFF and DEF sheets of the
# Target Depth Data workbook to reflect the newly simulated game.

$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 ("H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 408
$wsOff.Range("C2").Value = 315
$wsOff.Range("D2").Value = 103
$wsOff.Range("E2").Value = 55

# --- DEF sheet: row 2 ("H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 561
$wsDef.Range("C2").Value = 389
$wsDef.Range("D2").Value = 129
$wsDef.Range("E2").Value = 57
$wsDef.Range("F2").Value = 8
